$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.978.14'
$ws.Range("E2").Value = '  +3.41%  '
$ws.Range("D3").Value = '3.389.48'
$ws.Range("E3").Value = '  +3.16%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.53%  '
$ws.Range("D8").Value = '3.381.12'
$ws.Range("E8").Value = '  +3.09%  '
$ws.Range("E10").Value = '  +10.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.631'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.75%  '
$ws.Range("E13").Value = '  +6.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.70%  '
$ws.Range("D15").Value = '3.920.42'
$ws.Range("E15").Value = '  +3.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.120'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.26%  '
$ws.Range("D18").Value = '3.373.21'
$ws.Range("E18").Value = '  +2.82%  '
$ws.Range("D19").Value = '65.104.67'
$ws.Range("E19").Value = '  +3.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.994'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '468.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +13.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.89'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.55'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.69'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.50%  '
$ws.Range("E32").Value = '  +2.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '572.52'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '61.47'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.90%  '
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.79%  '
$ws.Range("E38").Value = '  -3.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.40%  '
$ws.Range("D40").Value = '0.0₃0745'
$ws.Range("E40").Value = '  +2.14%  '
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("D42").Value = '3.094.65'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.65%  '
$ws.Range("E45").Value = '  +4.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.135'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.01%  '
$ws.Range("E47").Value = '  +2.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.61'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '137.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.29'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.55%  '
